# ===========================================================================
# "se modifico un titulo" - TALLER SQL.docx edit
#
#  1) Collapse the spell-checked run fragments of the
#     "LEFT JOIN  productos on id_produ = id_productos" paragraph into a
#     single plain run (drop the <w:proofErr/> anchors).
#  2) Drop the stray _GoBack bookmark that sat in the "Nota:" paragraph.
#  3) Turn the last empty "List Paragraph" right above the foreign-key
#     question into a bold "PARTE 4:" heading.
#  4) Number the two remaining questions ("1." / "2."), splitting each into
#     two runs, and re-seat the _GoBack bookmark between "2." and the rest
#     of its question text.
# ===========================================================================

$d = $word.ActiveDocument

function Get-ParaOpenTag($para) {
    # Pull the paragraph's own <w:p .../> opening tag (paraId/textId/rsid...)
    # straight out of Word so the rebuilt paragraph keeps its identity.
    $owx = $para.Range.WordOpenXML
    $bodyIdx = $owx.IndexOf("<w:body>") + 8
    $closeIdx = $owx.IndexOf(">", $bodyIdx)
    return $owx.Substring($bodyIdx, $closeIdx - $bodyIdx + 1)
}

function Find-Paragraph($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Text -like $pattern) {
            return $cand
        }
    }
    return $null
}

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:xml="http://www.w3.org/XML/1998/namespace"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) Rebuild "LEFT JOIN  productos on id_produ = id_productos" as one run.
#    The range has to reach back into the previous paragraph's mark, or the
#    dangling <w:proofErr/> markers survive the rewrite.
# ---------------------------------------------------------------------------
$joinProdPara = Find-Paragraph "LEFT JOIN*productos on id_produ*id_productos*"
$prevPara = $joinProdPara.Previous()

$prevTag = Get-ParaOpenTag $prevPara
$curTag = Get-ParaOpenTag $joinProdPara

$prevBody = $prevTag + '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="1080"/><w:jc w:val="both"/></w:pPr>' + `
  '<w:r><w:t>LEFT JOIN  productos_categorias on id_cat = id_categorias</w:t></w:r></w:p>'
$curBody = $curTag + '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="1080"/><w:jc w:val="both"/></w:pPr>' + `
  '<w:r><w:t>LEFT JOIN  productos on id_produ = id_productos</w:t></w:r></w:p>'

$full = $d.Range($prevPara.Range.Start, $joinProdPara.Range.End)
$full.InsertXML($xmlHeader + $prevBody + $curBody + $xmlFooter)

# ---------------------------------------------------------------------------
# 2) Drop the stray _GoBack bookmark left in the "Nota:" paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) The last (empty) "List Paragraph" right before the foreign-key question
#    becomes a bold "PARTE 4:" heading.
# ---------------------------------------------------------------------------
$claveForaneaPara = Find-Paragraph "*claves for*neas?*"
$partePara = $claveForaneaPara.Previous()
$parteTag = Get-ParaOpenTag $partePara
$parteBody = $parteTag + '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="1080"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>PARTE 4:</w:t></w:r></w:p>'
$parteRange = $d.Range($partePara.Range.Start, $partePara.Range.End)
$parteRange.InsertXML($xmlHeader + $parteBody + $xmlFooter)

# ---------------------------------------------------------------------------
# 4) "4 ¿Para qué se usan las claves foráneas?" -> "1. ¿Para qué...?"
# ---------------------------------------------------------------------------
$claveForaneaPara = Find-Paragraph "*claves for*neas?*"
$claveTag = Get-ParaOpenTag $claveForaneaPara
$claveBody = $claveTag + '<w:pPr><w:jc w:val="both"/></w:pPr>' + `
  '<w:r><w:t>1.</w:t></w:r><w:r><w:t xml:space="preserve"> ¿Para qué se usan las claves foráneas?</w:t></w:r></w:p>'
$claveRange = $d.Range($claveForaneaPara.Range.Start, $claveForaneaPara.Range.End)
$claveRange.InsertXML($xmlHeader + $claveBody + $xmlFooter)

# ---------------------------------------------------------------------------
# 5) "5 ¿Qué diferencia hay entre INNER JOIN y LEFT JOIN?" -> "2. ¿Qué...?"
#    with the _GoBack bookmark re-seated between the "2." and the question.
# ---------------------------------------------------------------------------
$diffPara = Find-Paragraph "*diferencia hay entre INNER JOIN*"
$diffTag = Get-ParaOpenTag $diffPara
$diffBody = $diffTag + '<w:pPr><w:jc w:val="both"/></w:pPr>' + `
  '<w:r><w:t>2.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '<w:r><w:t xml:space="preserve"> ¿Qué diferencia hay entre INNER JOIN y LEFT JOIN?</w:t></w:r></w:p>'
$diffRange = $d.Range($diffPara.Range.Start, $diffPara.Range.End)
$diffRange.InsertXML($xmlHeader + $diffBody + $xmlFooter)

Write-Output "edit complete"
